# Apply the NASEM instructions workbook edits to the "Code sets" sheet:
#  - F17: "POTS" -> "pots"
#  - F23: "mood disorders" -> "mood disorder"
#  - F24: "inserstitial lung disease" -> "interstitial lung disease"
#  - Update the sheet view: scroll so row 2 is the top row, and move the
#    active selection to F23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Code sets")

$ws.Range("F17").Value = "pots"
$ws.Range("F24").Value = "interstitial lung disease"
$ws.Range("F23").Value = "mood disorder"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("F23").Select()
